# Update HRS yearly financials with latest figures ("Doing Updates for Financials")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Balance Sheet - Assets
$ws.Range("D43").Value = 1691000   # Net Receivables
$ws.Range("D44").Value = 942000    # Inventory
$ws.Range("D45").Value = 103000    # Other Current Assets
$ws.Range("D46").Value = 2224000   # Total Current Assets
$ws.Range("D49").Value = 7350000   # Goodwill
$ws.Range("D52").Value = 366000    # Other Assets
$ws.Range("D54").Value = 9851000   # Total Assets

# Balance Sheet - Liabilities
$ws.Range("D59").Value = 1160000   # Other Current Liabilities
$ws.Range("D60").Value = 1850000   # Total Current Liabilities
$ws.Range("D62").Value = 1315000   # Other Liabilities
$ws.Range("D66").Value = 6573000   # Total Liabilities

# Balance Sheet - Stockholders' Equity
$ws.Range("D72").Value = 1648000   # Retained Earnings
$ws.Range("D76").Value = 3278000   # Total Stockholder Equity

# Cash Flow Statement - Capital Expenditures
$ws.Range("I91").Value = -164800
$ws.Range("J91").Value = -209900
